# "Fruta / hortaliza, semanal"
# The weekly price list for this subset (Zapallo italiano, Mapocho Venta
# Directa de Santiago) was re-sorted: the 14 data rows (rows 2-15) keep the
# same set of records, but each record now lives on a different row.
# The common columns (A, B, C, E, F, G, H, I, R) are identical for every
# record, so the only columns that actually move are D (Fecha) and
# J:Q (Volumen .. Kg o Unidades).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot the current (pre-edit) values of the data rows (2-15), columns D and J:Q
$snapshot = @{}
$snapshot[2] = @{
  "D" = $ws.Range("D2").Value()
  "J" = $ws.Range("J2").Value()
  "K" = $ws.Range("K2").Value()
  "L" = $ws.Range("L2").Value()
  "M" = $ws.Range("M2").Value()
  "N" = $ws.Range("N2").Value()
  "O" = $ws.Range("O2").Value()
  "P" = $ws.Range("P2").Value()
  "Q" = $ws.Range("Q2").Value()
}
$snapshot[3] = @{
  "D" = $ws.Range("D3").Value()
  "J" = $ws.Range("J3").Value()
  "K" = $ws.Range("K3").Value()
  "L" = $ws.Range("L3").Value()
  "M" = $ws.Range("M3").Value()
  "N" = $ws.Range("N3").Value()
  "O" = $ws.Range("O3").Value()
  "P" = $ws.Range("P3").Value()
  "Q" = $ws.Range("Q3").Value()
}
$snapshot[4] = @{
  "D" = $ws.Range("D4").Value()
  "J" = $ws.Range("J4").Value()
  "K" = $ws.Range("K4").Value()
  "L" = $ws.Range("L4").Value()
  "M" = $ws.Range("M4").Value()
  "N" = $ws.Range("N4").Value()
  "O" = $ws.Range("O4").Value()
  "P" = $ws.Range("P4").Value()
  "Q" = $ws.Range("Q4").Value()
}
$snapshot[5] = @{
  "D" = $ws.Range("D5").Value()
  "J" = $ws.Range("J5").Value()
  "K" = $ws.Range("K5").Value()
  "L" = $ws.Range("L5").Value()
  "M" = $ws.Range("M5").Value()
  "N" = $ws.Range("N5").Value()
  "O" = $ws.Range("O5").Value()
  "P" = $ws.Range("P5").Value()
  "Q" = $ws.Range("Q5").Value()
}
$snapshot[6] = @{
  "D" = $ws.Range("D6").Value()
  "J" = $ws.Range("J6").Value()
  "K" = $ws.Range("K6").Value()
  "L" = $ws.Range("L6").Value()
  "M" = $ws.Range("M6").Value()
  "N" = $ws.Range("N6").Value()
  "O" = $ws.Range("O6").Value()
  "P" = $ws.Range("P6").Value()
  "Q" = $ws.Range("Q6").Value()
}
$snapshot[7] = @{
  "D" = $ws.Range("D7").Value()
  "J" = $ws.Range("J7").Value()
  "K" = $ws.Range("K7").Value()
  "L" = $ws.Range("L7").Value()
  "M" = $ws.Range("M7").Value()
  "N" = $ws.Range("N7").Value()
  "O" = $ws.Range("O7").Value()
  "P" = $ws.Range("P7").Value()
  "Q" = $ws.Range("Q7").Value()
}
$snapshot[8] = @{
  "D" = $ws.Range("D8").Value()
  "J" = $ws.Range("J8").Value()
  "K" = $ws.Range("K8").Value()
  "L" = $ws.Range("L8").Value()
  "M" = $ws.Range("M8").Value()
  "N" = $ws.Range("N8").Value()
  "O" = $ws.Range("O8").Value()
  "P" = $ws.Range("P8").Value()
  "Q" = $ws.Range("Q8").Value()
}
$snapshot[9] = @{
  "D" = $ws.Range("D9").Value()
  "J" = $ws.Range("J9").Value()
  "K" = $ws.Range("K9").Value()
  "L" = $ws.Range("L9").Value()
  "M" = $ws.Range("M9").Value()
  "N" = $ws.Range("N9").Value()
  "O" = $ws.Range("O9").Value()
  "P" = $ws.Range("P9").Value()
  "Q" = $ws.Range("Q9").Value()
}
$snapshot[10] = @{
  "D" = $ws.Range("D10").Value()
  "J" = $ws.Range("J10").Value()
  "K" = $ws.Range("K10").Value()
  "L" = $ws.Range("L10").Value()
  "M" = $ws.Range("M10").Value()
  "N" = $ws.Range("N10").Value()
  "O" = $ws.Range("O10").Value()
  "P" = $ws.Range("P10").Value()
  "Q" = $ws.Range("Q10").Value()
}
$snapshot[11] = @{
  "D" = $ws.Range("D11").Value()
  "J" = $ws.Range("J11").Value()
  "K" = $ws.Range("K11").Value()
  "L" = $ws.Range("L11").Value()
  "M" = $ws.Range("M11").Value()
  "N" = $ws.Range("N11").Value()
  "O" = $ws.Range("O11").Value()
  "P" = $ws.Range("P11").Value()
  "Q" = $ws.Range("Q11").Value()
}
$snapshot[12] = @{
  "D" = $ws.Range("D12").Value()
  "J" = $ws.Range("J12").Value()
  "K" = $ws.Range("K12").Value()
  "L" = $ws.Range("L12").Value()
  "M" = $ws.Range("M12").Value()
  "N" = $ws.Range("N12").Value()
  "O" = $ws.Range("O12").Value()
  "P" = $ws.Range("P12").Value()
  "Q" = $ws.Range("Q12").Value()
}
$snapshot[13] = @{
  "D" = $ws.Range("D13").Value()
  "J" = $ws.Range("J13").Value()
  "K" = $ws.Range("K13").Value()
  "L" = $ws.Range("L13").Value()
  "M" = $ws.Range("M13").Value()
  "N" = $ws.Range("N13").Value()
  "O" = $ws.Range("O13").Value()
  "P" = $ws.Range("P13").Value()
  "Q" = $ws.Range("Q13").Value()
}
$snapshot[14] = @{
  "D" = $ws.Range("D14").Value()
  "J" = $ws.Range("J14").Value()
  "K" = $ws.Range("K14").Value()
  "L" = $ws.Range("L14").Value()
  "M" = $ws.Range("M14").Value()
  "N" = $ws.Range("N14").Value()
  "O" = $ws.Range("O14").Value()
  "P" = $ws.Range("P14").Value()
  "Q" = $ws.Range("Q14").Value()
}
$snapshot[15] = @{
  "D" = $ws.Range("D15").Value()
  "J" = $ws.Range("J15").Value()
  "K" = $ws.Range("K15").Value()
  "L" = $ws.Range("L15").Value()
  "M" = $ws.Range("M15").Value()
  "N" = $ws.Range("N15").Value()
  "O" = $ws.Range("O15").Value()
  "P" = $ws.Range("P15").Value()
  "Q" = $ws.Range("Q15").Value()
}

# Re-order the rows: each destination row receives the values that
# previously belonged to the source row given by the mapping below.
$rowMap = @{
  2 = 9
  3 = 2
  4 = 5
  5 = 12
  6 = 11
  7 = 6
  8 = 10
  9 = 14
  10 = 15
  11 = 13
  12 = 3
  13 = 8
  14 = 7
  15 = 4
}

foreach ($dest in $rowMap.Keys) {
  $src = $rowMap[$dest]
  $row = $snapshot[$src]
  $ws.Range("D$dest").Value = $row["D"]
  $ws.Range("J$dest").Value = $row["J"]
  $ws.Range("K$dest").Value = $row["K"]
  $ws.Range("L$dest").Value = $row["L"]
  $ws.Range("M$dest").Value = $row["M"]
  $ws.Range("N$dest").Value = $row["N"]
  $ws.Range("O$dest").Value = $row["O"]
  $ws.Range("P$dest").Value = $row["P"]
  $ws.Range("Q$dest").Value = $row["Q"]
}